# Diary Guowei Li 01162020
# Insert a new diary entry row (Jan 16th, 2020) right after the first entry (row 10),
# within the templated table area (rows 11-39), shifting existing rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 11:39 down by one row (within the table's column span) to make room
# for the new entry, matching the same mechanism Excel uses for "Insert Copied Cells".
$ws.Range("A11:G39").Insert()

# Copy the formatting of the first data row (row 10) onto the newly inserted row 11.
$ws.Range("A10:G10").Copy()
$ws.Range("A11:G11").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new diary entry content.
$ws.Range("A11").Value = "Jan 16th, 2020"
$ws.Range("B11").Value = "5 - 8 pm"
$ws.Range("C11").Value = "N/A"
$ws.Range("D11").Value = "Revise last week's material, learn the basic strategies for code comprehension,  do an in-class practice and listen to a speech by Ping."
$ws.Range("E11").Value = "Understood different kinds of strategies for reading code, gained hands-on experience by doing practice, and also knew how professional programmers read code. "
$ws.Range("F11").Value = "First of all, different kinds of stategies can all be used. We don't need to insist on one specific strategy when reading code. Second, it's really necessary to assign meaningful names to variables, functions and classes. It helps a lot when people try to understand your code. Last but not least,  documenting is good for everybody in your group at work."
$ws.Range("G11").Value = "Great. I need time to get truly familiar with what I learned this time."

# Match the row height Excel computed for the new wrapped-text entry.
$ws.Rows("11:11").RowHeight = 117

# Update the view to reflect where the user ended up after editing.
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("G11").Select()
